# Updated cryptos list on Tue Sep  3 11:37:09 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.935.32"
$ws.Range("E2").Value = "  +0.99%  "

$ws.Range("D3").Value = "2.503.60"
$ws.Range("E3").Value = "  +0.65%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'534.71"
$ws.Range("E5").Value = "  +2.81%  "

$ws.Range("D6").Value = "'134.20"
$ws.Range("E6").Value = "  +1.58%  "

$ws.Range("E7").Value = "  +0.44%  "

$ws.Range("E8").Value = "  +2.41%  "

$ws.Range("D9").Value = "2.507.87"
$ws.Range("E9").Value = "  -0.64%  "

$ws.Range("D10").Value = "'0.0996"
$ws.Range("E10").Value = "  +1.95%  "

$ws.Range("E11").Value = "  -2.77%  "

$ws.Range("D12").Value = "'5.18"
$ws.Range("E12").Value = "  -0.82%  "

$ws.Range("E13").Value = "  -1.65%  "

$ws.Range("D14").Value = "2.947.85"
$ws.Range("E14").Value = "  +0.21%  "

$ws.Range("D15").Value = "58.741.93"
$ws.Range("E15").Value = "  +0.84%  "

$ws.Range("D16").Value = "'22.39"
$ws.Range("E16").Value = "  +0.49%  "

$ws.Range("E17").Value = "  +0.29%  "

$ws.Range("D18").Value = "2.503.71"
$ws.Range("E18").Value = "  -0.40%  "

$ws.Range("D19").Value = "'10.64"
$ws.Range("E19").Value = "  -1.13%  "

$ws.Range("D20").Value = "'4.26"
$ws.Range("E20").Value = "  +1.41%  "

$ws.Range("D21").Value = "'321.24"
$ws.Range("E21").Value = "  -1.11%  "

$ws.Range("D22").Value = "'6.22"
$ws.Range("E22").Value = "  +2.42%  "

$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  +0.26%  "

$ws.Range("D24").Value = "'65.85"
$ws.Range("E24").Value = "  +3.26%  "

$ws.Range("E25").Value = "  +0.54%  "

$ws.Range("E26").Value = "  +1.38%  "

$ws.Range("E27").Value = "  -1.43%  "

$ws.Range("D28").Value = "'7.46"
$ws.Range("E28").Value = "  +0.88%  "

$ws.Range("D29").Value = "0.0₃0756"
$ws.Range("E29").Value = "  +0.86%  "

$ws.Range("D30").Value = "'172.05"
$ws.Range("E30").Value = "  +2.37%  "

$ws.Range("E31").Value = "  +1.42%  "

$ws.Range("D32").Value = "'6.29"
$ws.Range("E32").Value = "  -0.03%  "

$ws.Range("E33").Value = "  -0.67%  "

$ws.Range("E34").Value = "  +0.05%  "

$ws.Range("E35").Value = "  +0.40%  "

$ws.Range("E36").Value = "  +0.14%  "

$ws.Range("E37").Value = "  -3.80%  "

$ws.Range("E38").Value = "  -0.02%  "

$ws.Range("E39").Value = "  +3.41%  "

$ws.Range("E40").Value = "  +5.77%  "

$ws.Range("D41").Value = "'36.53"
$ws.Range("E41").Value = "  -0.45%  "

$ws.Range("E42").Value = "  +0.95%  "

$ws.Range("D43").Value = "'274.89"
$ws.Range("E43").Value = "  -1.69%  "

$ws.Range("D44").Value = "'131.05"
$ws.Range("E44").Value = "  +6.99%  "

$ws.Range("E45").Value = "  -1.96%  "

$ws.Range("E46").Value = "  -1.49%  "

$ws.Range("D47").Value = "'0.0937"
$ws.Range("E47").Value = "  +1.65%  "

$ws.Range("D48").Value = "'0.0511"
$ws.Range("E48").Value = "  +2.14%  "

$ws.Range("D49").Value = "'0.0218"
$ws.Range("E49").Value = "  +1.98%  "

$ws.Range("D50").Value = "'16.81"
$ws.Range("E50").Value = "  -1.56%  "

$ws.Range("D51").Value = "1.749.01"
$ws.Range("E51").Value = "  +0.11%  "
